$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, border, centered) from an existing header cell
$ws.Range("AA1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Fill in season record values (Wins, Losses, Ties) for every data row
for ($row = 2; $row -le 52; $row++) {
    $ws.Cells.Item($row, 30).Value = 93
    $ws.Cells.Item($row, 31).Value = 69
    $ws.Cells.Item($row, 32).Value = 0
}
